$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two columns: A = Ticker Symbol (header + 12 tickers),
# B = 2014 (header + 12 negative EPS values). The ticker-symbol list is no
# longer needed, so drop column A entirely; column B (the "2014" header and
# its EPS figures) shifts left and becomes the new column A.
$ws.Range("A:A").Delete()

# Call out the worst (most negative) EPS figure in the remaining data with a
# red fill. That's -14.06, which landed on row 3 after the shift.
$ws.Range("A3").Interior.Color = 255
